$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows for 2021-10-07 and 2021-10-08 (the "recup"/"verlof"
# placeholder days). Deleting entire rows shifts everything below them up,
# so the remaining 2021-11-11 entry moves from row 9 to row 7.
$ws.Range("A7:A8").EntireRow.Delete()

# Remove the now-unused "verlof" and "recup" columns (E and F) along with
# their header cells and shared-string entries.
$ws.Range("E1:F1").EntireColumn.Delete()

# Match the workbook's saved selection state.
$ws.Range("D8").Select()
